$d = $word.ActiveDocument

# --------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the first
#    (Heading1) paragraph.
# --------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"
$metaRange = $metaPara.Range

# Fill the paragraph with the full text first (plain formatting) ...
$insertionPoint = $d.Range($metaRange.End - 1, $metaRange.End - 1)
$insertionPoint.Text = "Meta description: Join treasure hunters in Book of Ra Temple of Gold, a Novomatic slot game with high volatility and expanding symbols. Play for free and win big!"

# ... then bold just the "Meta description" label.
$labelLength = ([string]"Meta description").Length
$labelRange = $d.Range($metaRange.Start, $metaRange.Start + $labelLength)
$labelRange.Font.Bold = 1

# --------------------------------------------------------------------
# 2) Remove the duplicated bold "Play Book of Ra Temple of Gold..."
#    paragraph that used to sit at the bottom of the document.
# --------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
  $p = $d.Paragraphs($i)
  if ($p.Range.Text.TrimEnd() -eq "Play Book of Ra Temple of Gold Free - Exciting Adventures and High Returns" -and $i -ne 1) {
    $p.Range.Delete()
    break
  }
}

# --------------------------------------------------------------------
# 3) Replace the text of the closing italic paragraph with the new
#    image-generation prompt, keeping its italic formatting.
# --------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
  $p = $d.Paragraphs($i)
  if ($p.Range.Text.TrimEnd() -eq "Join treasure hunters in Book of Ra Temple of Gold, a Novomatic slot game with high volatility and expanding symbols. Play for free and win big!") {
    $targetRange = $p.Range
    [void]$targetRange.MoveEnd(1, -1)  # exclude the paragraph mark
    $targetRange.Text = "Prompt: Create a feature image for Book of Ra: Temple of Gold that captures the adventurous spirit of the game while showcasing the happy Maya warrior with glasses in a cartoon style. The feature image should include a jungle background with hints of gold and a temple towering in the distance. The foreground should have the happy Maya warrior holding a treasure chest filled with gold coins while wearing a backpack and a pair of glasses. The warrior should be positioned in a triumphant stance, with one arm raised in celebration. The warrior's clothing should consist of traditional Maya attire with a modern twist as depicted in the game. The image should be in cartoon style, with bold colors and shading to enhance the visual impact of the image."
    break
  }
}

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
